# Slide 45 ("Principles"), Content Placeholder:
#   - split the existing "SOLID Principles" run (which carries the
#     rId3 hyperlink) into two runs: "SOLID " and "Principles"
#   - append new bullets after it:
#       Object-Oriented Principles
#         Objects, classes        (indented)
#         Encapsulation           (indented)
#         Inheritance             (indented)
#         Polymorphism            (indented)
#         Design Patterns         (indented)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(45)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 7"
$tf = $sh.TextFrame
$tr = $tf.TextRange

$hyperlinkUrl = "https://en.wikipedia.org/wiki/SOLID_(object-oriented_design)"

# --- Remove the "SOLID Principles" paragraph (including its paragraph
#     mark) so the new bullets can be typed in right after "DRY ..."
#     without inheriting the SOLID-Principles hyperlink formatting. ---
$para1 = $tr.Paragraphs(1, 1)
$para2 = $tr.Paragraphs(2, 1)
$solidRun = $tr.Characters($para1.Start + $para1.Length, $para2.Length + 1)
$solidRun.Delete()

# --- Re-type the six new, plain (non-hyperlinked) paragraphs right
#     after the "DRY ..." paragraph. ---
$para1 = $tr.Paragraphs(1, 1)
$para1.InsertAfter("`rObject-Oriented Principles`rObjects, classes`rEncapsulation`rInheritance`rPolymorphism`rDesign Patterns")

# --- Re-type "SOLID Principles" as its own paragraph, right after
#     "DRY ..." and before the new bullets, still with no hyperlink. ---
$para1 = $tr.Paragraphs(1, 1)
$para1.InsertAfter("`rSOLID Principles")

# --- Re-apply the hyperlink to the whole "SOLID Principles" run; this
#     reuses the presentation's existing rId3 relationship (same URL)
#     instead of creating a new one. ---
$para2 = $tr.Paragraphs(2, 1)
$para2.ActionSettings(1).Hyperlink.Address = $hyperlinkUrl

# --- Split "SOLID Principles" into "SOLID " + "Principles", two runs,
#     both keeping the rId3 hyperlink. ---
$para2 = $tr.Paragraphs(2, 1)
$firstPart = $tr.Characters($para2.Start, 6)
$firstPart.Text = "SOLID "

# --- Demote the five sub-bullets under "Object-Oriented Principles". ---
for ($i = 4; $i -le 8; $i++) {
    $tr.Paragraphs($i, 1).IndentLevel = 2
}

Write-Output $tr.Text
